$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.141279578208923
$ws.Range("B1").Value = 1.77022397518158
$ws.Range("C1").Value = 3.977267980575562
$ws.Range("D1").Value = 2.755751132965088
$ws.Range("E1").Value = 0.2947202026844025
